$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

$xlShiftToRight = [Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight

# ---------------------------------------------------------------------
# 1) "target" list (column A): insert a new category entry "text" right
#    after "step" (A24), pushing the remaining entries (web, webalert,
#    webcookie, ws, ws.async, xml -> A25..A30) down by one row (A26..A31).
#    Cell-by-cell copy (bottom-up) so that the other, unrelated columns on
#    those same rows are left completely untouched.
# ---------------------------------------------------------------------
for ($r = 30; $r -ge 25; $r--) {
    $ws.Cells.Item($r + 1, 1).Value2 = $ws.Cells.Item($r, 1).Value2
}
$ws.Cells.Item(25, 1).Value2 = "text"

# ---------------------------------------------------------------------
# 2) "base" list (column E): insert a new function entry
#    "outputToCloud(resource)" right after "macro(file,sheet,name)" (E21),
#    pushing the remaining entries (E22..E38) down by one row (E23..E39).
# ---------------------------------------------------------------------
for ($r = 38; $r -ge 22; $r--) {
    $ws.Cells.Item($r + 1, 5).Value2 = $ws.Cells.Item($r, 5).Value2
}
$ws.Cells.Item(22, 5).Value2 = "outputToCloud(resource)"

# ---------------------------------------------------------------------
# 3) Insert a brand-new column at Y to host the "text" function-family
#    table; this pushes the existing web/webalert/webcookie/ws/ws.async/xml
#    columns one column to the right (Y->Z, Z->AA, ... AD->AE). A whole
#    column insert is naturally scoped to columns Y.. onward, so it does
#    not disturb columns A..X.
# ---------------------------------------------------------------------
$ws.Range("Y1").EntireColumn.Insert($xlShiftToRight)
$ws.Cells.Item(1, 25).Value2 = "text"
$ws.Cells.Item(2, 25).Value2 = "spellCheck(var,profile,text)"

# ---------------------------------------------------------------------
# 4) Update the named ranges so they reflect the new row/column positions.
# ---------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$E`$2:`$E`$39"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$31"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AC`$2:`$AC`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AD`$2:`$AD`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AE`$2:`$AE`$27"
$wb.Names.Add("text", "='#system'!`$Y`$2:`$Y`$2")
